# "+ Ban nop/Phieu van dap.xlsx" - cap nhat bang diem cham van dap.
# Chi co 2 o nhap lieu thuc su thay doi (muc 4.4 -> 4.2 cham lai diem danh gia):
#   D45 : 1    -> 0.75
#   D50 : 0.5  -> 0.1
# Cac cong thuc phu thuoc (G42, G46, G11, H11) tu dong tinh lai.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhieuVanDap")

$ws.Range("D45").Value = 0.75
$ws.Range("D50").Value = 0.1

# Dua view ve dau trang (bo trang thai cuon/o dang chon cu tai H44)
# truoc khi luu ban nop chinh thuc.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
[void]$ws.Range("A1").Select()
